$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused / Clear Glass Lens
$ws.Range("H33").Value = 552.5
$ws.Range("J33").Value = 1169
$ws.Range("L33").Value = 1169
$ws.Range("N33").Value = -1627

# Row 40: Stuck in the Moment / Horn Glue
$ws.Range("H40").Value = 4998
$ws.Range("J40").Value = 4998.5
$ws.Range("L40").Value = 4998.5
$ws.Range("N40").Value = -5348.5

# Row 53: No Accounting for Waste / Enchanted Electrum Ink
$ws.Range("H53").Value = 4741
$ws.Range("I53").Value = 6973.5
$ws.Range("J53").Value = 276
$ws.Range("K53").Value = 6973.5
$ws.Range("L53").Value = 276
$ws.Range("M53").Value = -6336.5
$ws.Range("N53").Value = -1550

# Row 93: Spellbound / Koppranickel Index
$ws.Range("H93").Value = 120219900
$ws.Range("J93").Value = 120219900
$ws.Range("L93").Value = 120219900
$ws.Range("N93").Value = -120224892

# Row 94: Magic Beans / Growth Formula Eta
$ws.Range("H94").Value = 62656212
$ws.Range("I94").Value = 100003940
$ws.Range("J94").Value = 410002
$ws.Range("K94").Value = 100003940
$ws.Range("L94").Value = 410002
$ws.Range("M94").Value = -100003489
$ws.Range("N94").Value = -410904

# Row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 1445.6
$ws.Range("J112").Value = 1497.9231
$ws.Range("L112").Value = 4493.7693
$ws.Range("N112").Value = -6709.7693

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 5725.3335
$ws.Range("I132").Value = 5114.7617
$ws.Range("K132").Value = 15344.2851
$ws.Range("M132").Value = -12814.2851

# Row 140: Tome for Tradition / Book of Ra'Kaznar
$ws.Range("H140").Value = 80697.39999999999
$ws.Range("J140").Value = 80697.39999999999
$ws.Range("L140").Value = 80697.39999999999
$ws.Range("N140").Value = -91057.39999999999

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 2490.2778
$ws.Range("I32").Value = 2519.1177
$ws.Range("K32").Value = 2519.1177
$ws.Range("M32").Value = -2232.1177

# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 7649.5454
$ws.Range("J45").Value = 7416.1665
$ws.Range("L45").Value = 7416.1665
$ws.Range("N45").Value = -8170.1665

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 4977
$ws.Range("I61").Value = 4388.3184
$ws.Range("J61").Value = 5973.231
$ws.Range("K61").Value = 4388.3184
$ws.Range("L61").Value = 5973.231
$ws.Range("M61").Value = -4176.3184
$ws.Range("N61").Value = -6397.231

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 4977
$ws.Range("I136").Value = 4388.3184
$ws.Range("J136").Value = 5973.231
$ws.Range("K136").Value = 13164.9552
$ws.Range("L136").Value = 17919.693
$ws.Range("M136").Value = -10614.9552
$ws.Range("N136").Value = -23019.693

# Row 139: Backing up My Words / Titanium Gold Thornplate of Fending
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 82: Spirituality Inspector / Titanium Lump Hammer
$ws.Range("H82").Value = 62268.145
$ws.Range("J82").Value = 75399.39999999999
$ws.Range("L82").Value = 75399.39999999999
$ws.Range("N82").Value = -76165.39999999999

# Row 85: The Clamor for Hammers (L) / Titanium Lump Hammer
$ws.Range("H85").Value = 62268.145
$ws.Range("J85").Value = 75399.39999999999
$ws.Range("L85").Value = 75399.39999999999
$ws.Range("N85").Value = -78051.39999999999

# Row 100: And My Axe / Doman Iron War Axe
$ws.Range("H100").Value = 38475
$ws.Range("J100").Value = 38475
$ws.Range("L100").Value = 38475
$ws.Range("N100").Value = -40639

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 14748.5
$ws.Range("I107").Value = 15426.857
$ws.Range("K107").Value = 15426.857
$ws.Range("M107").Value = -13506.857

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 1656.9348
$ws.Range("I134").Value = 1174.2162
$ws.Range("K134").Value = 3522.6486
$ws.Range("M134").Value = -987.6486000000004

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 3039.9473
$ws.Range("J31").Value = 5172.6665
$ws.Range("L31").Value = 5172.6665
$ws.Range("N31").Value = -5762.6665

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 3039.9473
$ws.Range("J34").Value = 5172.6665
$ws.Range("L34").Value = 5172.6665
$ws.Range("N34").Value = -5576.6665

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 53777.777
$ws.Range("I132").Value = 4800
$ws.Range("K132").Value = 14400
$ws.Range("M132").Value = -11870

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 7396.467
$ws.Range("I134").Value = 7326.1113
$ws.Range("K134").Value = 21978.3339
$ws.Range("M134").Value = -19443.3339

$ws = $wb.Worksheets.Item("GSM")
# Row 5: Hora at Me / Bone Hora
$ws.Range("H5").Value = 6000
$ws.Range("I5").Value = 6000
$ws.Range("K5").Value = 6000
$ws.Range("M5").Value = -5888

# Row 57: Gold Is So Last Year / Electrum Circlet (Amber)
$ws.Range("H57").Value = 28995.8
$ws.Range("J57").Value = 34993
$ws.Range("L57").Value = 34993
$ws.Range("N57").Value = -36633

# Row 59: Sew Not Doing This / Electrum Needle
$ws.Range("H59").Value = 8000
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 5824.478
$ws.Range("I70").Value = 5550.4287
$ws.Range("J70").Value = 6250.778
$ws.Range("K70").Value = 5550.4287
$ws.Range("L70").Value = 6250.778
$ws.Range("M70").Value = -5280.4287
$ws.Range("N70").Value = -6790.778

# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 5824.478
$ws.Range("I73").Value = 5550.4287
$ws.Range("J73").Value = 6250.778
$ws.Range("K73").Value = 5550.4287
$ws.Range("L73").Value = 6250.778
$ws.Range("M73").Value = -4614.4287
$ws.Range("N73").Value = -8122.778

# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 7495.087
$ws.Range("I102").Value = 10325.8
$ws.Range("J102").Value = 2187.5
$ws.Range("K102").Value = 10325.8
$ws.Range("L102").Value = 2187.5
$ws.Range("M102").Value = -8703.799999999999
$ws.Range("N102").Value = -5431.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 3370.5334
$ws.Range("I22").Value = 3478.3
$ws.Range("J22").Value = 3155
$ws.Range("K22").Value = 3478.3
$ws.Range("L22").Value = 3155
$ws.Range("M22").Value = -3183.3
$ws.Range("N22").Value = -3745

# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 3370.5334
$ws.Range("I27").Value = 3478.3
$ws.Range("J27").Value = 3155
$ws.Range("K27").Value = 3478.3
$ws.Range("L27").Value = 3155
$ws.Range("M27").Value = -3371.3
$ws.Range("N27").Value = -3369

# Row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 3822.3845
$ws.Range("J46").Value = 5624
$ws.Range("L46").Value = 5624
$ws.Range("N46").Value = -6000

# Row 105: Thick and Thin / Gazelleskin Corselet of Scouting
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 652123.2
$ws.Range("I132").Value = 1492264.4
$ws.Range("K132").Value = 4476793.199999999
$ws.Range("M132").Value = -4474263.199999999

$ws = $wb.Worksheets.Item("WVR")
# Row 44: Edmelle's Hair / Linen Wedge Cap of Gathering
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

# Row 47: The Wages of Sin / Linen Coatee of Crafting
$ws.Range("H47").Value = 20000
$ws.Range("J47").Value = 20000
$ws.Range("L47").Value = 20000
$ws.Range("N47").Value = -21144

# Row 62: Pride Up in Smoke / Rainbow Cloth
$ws.Range("H62").Value = 273300.2
$ws.Range("I62").Value = 904000.7
$ws.Range("K62").Value = 904000.7
$ws.Range("M62").Value = -903376.7

# Row 65: Desperate for Diversionaries (L) / Rainbow Cloth
$ws.Range("H65").Value = 273300.2
$ws.Range("I65").Value = 904000.7
$ws.Range("K65").Value = 4520003.5
$ws.Range("M65").Value = -4516883.5

# Row 92: Modest Beginnings / Bloodhempen Culottes of Casting
$ws.Range("H92").Value = 137537500
$ws.Range("J92").Value = 137537500
$ws.Range("L92").Value = 137537500
$ws.Range("N92").Value = -137542492

# Row 96: Skills on Display / Ruby Cotton Cloth
$ws.Range("H96").Value = 2092.1482
$ws.Range("I96").Value = 1754.6666
$ws.Range("J96").Value = 2767.111
$ws.Range("K96").Value = 1754.6666
$ws.Range("L96").Value = 2767.111
$ws.Range("M96").Value = -381.6666
$ws.Range("N96").Value = -5513.111

# Row 126: A Polished Purchase / Snow Linen
$ws.Range("H126").Value = 35929.848
$ws.Range("I126").Value = 47455.332
$ws.Range("J126").Value = 9997.5
$ws.Range("K126").Value = 142365.996
$ws.Range("L126").Value = 29992.5
$ws.Range("M126").Value = -139895.996
$ws.Range("N126").Value = -34932.5
